$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1025.4546
$ws.Range("I15").Value = 1025.4546
$ws.Range("K15").Value = 3076.3638
$ws.Range("M15").Value = -2907.3638
$ws.Range("H18").Value = 7102.067
$ws.Range("I18").Value = 466.35715
$ws.Range("K18").Value = 466.35715
$ws.Range("M18").Value = -182.35715
$ws.Range("H92").Value = 456.5
$ws.Range("I92").Value = 456.5
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 456.5
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 791.5
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 102230.4
$ws.Range("J113").Value = 2166.5
$ws.Range("L113").Value = 2166.5
$ws.Range("N113").Value = -8674.5
$ws.Range("H124").Value = 46551
$ws.Range("J124").Value = 46551
$ws.Range("L124").Value = 46551
$ws.Range("N124").Value = -56371
$ws.Range("H128").Value = 46248.285
$ws.Range("J128").Value = 46248.285
$ws.Range("L128").Value = 46248.285
$ws.Range("N128").Value = -56208.285
$ws.Range("H138").Value = 4284.909
$ws.Range("I138").Value = 2079.4375
$ws.Range("J138").Value = 5189.718
$ws.Range("K138").Value = 6238.3125
$ws.Range("L138").Value = 15569.154
$ws.Range("M138").Value = -1098.3125
$ws.Range("N138").Value = -25849.154
$ws.Range("H141").Value = 3484.762
$ws.Range("I141").Value = 3265.5557
$ws.Range("J141").Value = 4800
$ws.Range("K141").Value = 9796.667099999999
$ws.Range("L141").Value = 14400
$ws.Range("M141").Value = -4616.667099999999
$ws.Range("N141").Value = -24760

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 42935
$ws.Range("I97").Value = 50876
$ws.Range("J97").Value = 3230
$ws.Range("K97").Value = 50876
$ws.Range("L97").Value = 3230
$ws.Range("M97").Value = -50380
$ws.Range("N97").Value = -4222

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 9098.666999999999
$ws.Range("I122").Value = 9226
$ws.Range("J122").Value = 8971.333000000001
$ws.Range("K122").Value = 27678
$ws.Range("L122").Value = 26913.999
$ws.Range("M122").Value = -25228
$ws.Range("N122").Value = -31813.999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1042.7142
$ws.Range("I58").Value = 968.3333
$ws.Range("J58").Value = 1098.5
$ws.Range("K58").Value = 2904.9999
$ws.Range("L58").Value = 3295.5
$ws.Range("M58").Value = -2776.9999
$ws.Range("N58").Value = -3551.5
$ws.Range("H64").Value = 2500
$ws.Range("J64").Value = 2500
$ws.Range("L64").Value = 7500
$ws.Range("N64").Value = -8040
$ws.Range("H67").Value = 2500
$ws.Range("J67").Value = 2500
$ws.Range("L67").Value = 7500
$ws.Range("N67").Value = -9372
$ws.Range("H94").Value = 4388.5557
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 4785.2856
$ws.Range("K94").Value = 9000
$ws.Range("L94").Value = 14355.8568
$ws.Range("M94").Value = -8324
$ws.Range("N94").Value = -15707.8568
$ws.Range("H114").Value = 237.25
$ws.Range("I114").Value = 245.33333
$ws.Range("J114").Value = 213
$ws.Range("K114").Value = 735.99999
$ws.Range("L114").Value = 639
$ws.Range("M114").Value = 2518.00001
$ws.Range("N114").Value = -7147
$ws.Range("H131").Value = 826.46
$ws.Range("J131").Value = 837.43616
$ws.Range("L131").Value = 2512.30848
$ws.Range("N131").Value = -12592.30848
$ws.Range("H136").Value = 2280.2
$ws.Range("J136").Value = 3700
$ws.Range("L136").Value = 11100
$ws.Range("N136").Value = -21300

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 12581.6
$ws.Range("J48").Value = 12581.6
$ws.Range("L48").Value = 12581.6
$ws.Range("N48").Value = -13551.6
$ws.Range("H57").Value = 19800
$ws.Range("J57").Value = 19800
$ws.Range("L57").Value = 19800
$ws.Range("N57").Value = -21440
$ws.Range("H122").Value = 5161.8184
$ws.Range("I122").Value = 4357.143
$ws.Range("J122").Value = 6570
$ws.Range("K122").Value = 13071.429
$ws.Range("L122").Value = 19710
$ws.Range("M122").Value = -10621.429
$ws.Range("N122").Value = -24610

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 586.06665
$ws.Range("I22").Value = 449.5
$ws.Range("J22").Value = 607.0769
$ws.Range("K22").Value = 449.5
$ws.Range("L22").Value = 607.0769
$ws.Range("M22").Value = -154.5
$ws.Range("N22").Value = -1197.0769
$ws.Range("H27").Value = 586.06665
$ws.Range("I27").Value = 449.5
$ws.Range("J27").Value = 607.0769
$ws.Range("K27").Value = 449.5
$ws.Range("L27").Value = 607.0769
$ws.Range("M27").Value = -342.5
$ws.Range("N27").Value = -821.0769
$ws.Range("H46").Value = 1125412.2
$ws.Range("I46").Value = 650
$ws.Range("J46").Value = 1446772.9
$ws.Range("K46").Value = 650
$ws.Range("L46").Value = 1446772.9
$ws.Range("M46").Value = -462
$ws.Range("N46").Value = -1447148.9
$ws.Range("H74").Value = 19499.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 19499.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 19499.4
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -21495.4
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 19499.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 19499.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 58498.2
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -68482.20000000001
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H128").Value = 47890
$ws.Range("J128").Value = 47890
$ws.Range("L128").Value = 47890
$ws.Range("N128").Value = -57850

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 33141.8
$ws.Range("J117").Value = 33141.8
$ws.Range("L117").Value = 33141.8
$ws.Range("N117").Value = -42319.8
